$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 534.5
$ws.Range("I8").Value = 534.5
$ws.Range("K8").Value = 1603.5
$ws.Range("M8").Value = -1464.5
$ws.Range("H9").Value = 89.25
$ws.Range("I9").Value = 88.5
$ws.Range("K9").Value = 88.5
$ws.Range("M9").Value = 80.5
$ws.Range("H43").Value = 2399.6667
$ws.Range("I43").Value = 2799.5
$ws.Range("J43").Value = 1600
$ws.Range("K43").Value = 2799.5
$ws.Range("L43").Value = 1600
$ws.Range("M43").Value = -2730.5
$ws.Range("N43").Value = -1738
$ws.Range("H51").Value = 34999
$ws.Range("J51").Value = 34999
$ws.Range("L51").Value = 34999
$ws.Range("N51").Value = -35967
$ws.Range("H53").Value = 59.75
$ws.Range("J53").Value = 65
$ws.Range("L53").Value = 65
$ws.Range("N53").Value = -1339
$ws.Range("H80").Value = 724.6667
$ws.Range("J80").Value = 835.5714
$ws.Range("L80").Value = 2506.7142
$ws.Range("N80").Value = -4502.7142
$ws.Range("H83").Value = 724.6667
$ws.Range("J83").Value = 835.5714
$ws.Range("L83").Value = 7520.1426
$ws.Range("N83").Value = -17504.1426
$ws.Range("H94").Value = 2000
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("H107").Value = 171.6923
$ws.Range("I107").Value = 154
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 154
$ws.Range("L107").Value = 200
$ws.Range("M107").Value = 1766
$ws.Range("N107").Value = -4040
$ws.Range("H111").Value = 612.3333
$ws.Range("I111").Value = 612.3333
$ws.Range("K111").Value = 1836.9999
$ws.Range("M111").Value = 1230.0001
$ws.Range("H113").Value = 3332.5
$ws.Range("I113").Value = 3332
$ws.Range("J113").Value = 3333
$ws.Range("K113").Value = 3332
$ws.Range("L113").Value = 3333
$ws.Range("M113").Value = -78
$ws.Range("N113").Value = -9841
$ws.Range("H132").Value = 883.4
$ws.Range("I132").Value = 604.25
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 1812.75
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = 717.25
$ws.Range("N132").Value = -11060
$ws.Range("H138").Value = 4186.75
$ws.Range("J138").Value = 4642
$ws.Range("L138").Value = 13926
$ws.Range("N138").Value = -24206

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2657
$ws.Range("I61").Value = 1355.3334
$ws.Range("K61").Value = 1355.3334
$ws.Range("M61").Value = -1143.3334
$ws.Range("H74").Value = 1005.5
$ws.Range("I74").Value = 1012
$ws.Range("J74").Value = 999
$ws.Range("K74").Value = 1012
$ws.Range("L74").Value = 999
$ws.Range("M74").Value = -138
$ws.Range("N74").Value = -2747
$ws.Range("H77").Value = 1005.5
$ws.Range("I77").Value = 1012
$ws.Range("J77").Value = 999
$ws.Range("K77").Value = 5060
$ws.Range("L77").Value = 4995
$ws.Range("M77").Value = -692
$ws.Range("N77").Value = -13731
$ws.Range("H132").Value = 3849.9285
$ws.Range("I132").Value = 2883.5
$ws.Range("J132").Value = 4574.75
$ws.Range("K132").Value = 8650.5
$ws.Range("L132").Value = 13724.25
$ws.Range("M132").Value = -6120.5
$ws.Range("N132").Value = -18784.25
$ws.Range("H136").Value = 2657
$ws.Range("I136").Value = 1355.3334
$ws.Range("K136").Value = 4066.0002
$ws.Range("M136").Value = -1516.0002

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1045.75
$ws.Range("I99").Value = 1045.75
$ws.Range("K99").Value = 1045.75
$ws.Range("M99").Value = 452.25
$ws.Range("H107").Value = 33335404
$ws.Range("I107").Value = 1954.909
$ws.Range("J107").Value = 125002380
$ws.Range("K107").Value = 1954.909
$ws.Range("L107").Value = 125002380
$ws.Range("M107").Value = -34.90900000000011
$ws.Range("N107").Value = -125006220
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H134").Value = 10749.917
$ws.Range("I134").Value = 1499.8334
$ws.Range("K134").Value = 4499.5002
$ws.Range("M134").Value = -1964.5002

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1466.3334
$ws.Range("I2").Value = 1466.3334
$ws.Range("K2").Value = 1466.3334
$ws.Range("M2").Value = -1353.3334
$ws.Range("H17").Value = 649.5
$ws.Range("I17").Value = 649.5
$ws.Range("K17").Value = 649.5
$ws.Range("M17").Value = -475.5
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H99").Value = 2775.4167
$ws.Range("I99").Value = 2610.2222
$ws.Range("J99").Value = 3271
$ws.Range("K99").Value = 2610.2222
$ws.Range("L99").Value = 3271
$ws.Range("M99").Value = -1112.2222
$ws.Range("N99").Value = -6267
$ws.Range("H103").Value = 75999.5
$ws.Range("I103").Value = 57499.5
$ws.Range("J103").Value = 94499.5
$ws.Range("K103").Value = 57499.5
$ws.Range("L103").Value = 94499.5
$ws.Range("M103").Value = -56327.5
$ws.Range("N103").Value = -96843.5
$ws.Range("H126").Value = 2775.4167
$ws.Range("I126").Value = 2610.2222
$ws.Range("J126").Value = 3271
$ws.Range("K126").Value = 7830.6666
$ws.Range("L126").Value = 9813
$ws.Range("M126").Value = -5360.6666
$ws.Range("N126").Value = -14753
$ws.Range("H134").Value = 2030.5
$ws.Range("I134").Value = 2030.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6091.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3556.5
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5525792
$ws.Range("I4").Value = 643428.4399999999
$ws.Range("J4").Value = 11221883
$ws.Range("K4").Value = 1930285.32
$ws.Range("L4").Value = 33665649
$ws.Range("M4").Value = -1930173.32
$ws.Range("N4").Value = -33665873
$ws.Range("H5").Value = 1250
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 4500
$ws.Range("M5").Value = -2888
$ws.Range("N5").Value = -4724
$ws.Range("H15").Value = 320.8
$ws.Range("H29").Value = 86.5
$ws.Range("J29").Value = 79
$ws.Range("L29").Value = 237
$ws.Range("N29").Value = -791
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H99").Value = 4950
$ws.Range("I99").Value = 4616.6665
$ws.Range("K99").Value = 13849.9995
$ws.Range("M99").Value = -11603.9995
$ws.Range("H135").Value = 1250
$ws.Range("I135").Value = 1000
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 9000
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -6465
$ws.Range("N135").Value = -18570

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 550
$ws.Range("I5").Value = 550
$ws.Range("K5").Value = 550
$ws.Range("M5").Value = -438
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H107").Value = 1548
$ws.Range("I107").Value = 1548
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1548
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 372
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 4498.8335
$ws.Range("I132").Value = 4331.6665
$ws.Range("K132").Value = 12994.9995
$ws.Range("M132").Value = -10464.9995

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 694.7273
$ws.Range("I55").Value = 663.375
$ws.Range("J55").Value = 778.3333
$ws.Range("K55").Value = 663.375
$ws.Range("L55").Value = 778.3333
$ws.Range("M55").Value = -490.375
$ws.Range("N55").Value = -1124.3333
$ws.Range("H82").Value = 24131.666
$ws.Range("I82").Value = 13397.75
$ws.Range("K82").Value = 13397.75
$ws.Range("M82").Value = -13036.75
$ws.Range("H85").Value = 24131.666
$ws.Range("I85").Value = 13397.75
$ws.Range("K85").Value = 13397.75
$ws.Range("M85").Value = -12149.75
$ws.Range("H136").Value = 1974.7778
$ws.Range("I136").Value = 1973.3529
$ws.Range("J136").Value = 1999
$ws.Range("K136").Value = 5920.0587
$ws.Range("L136").Value = 5997
$ws.Range("M136").Value = -3370.0587
$ws.Range("N136").Value = -11097

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 50010164
$ws.Range("J94").Value = 50010164
$ws.Range("L94").Value = 50010164
$ws.Range("N94").Value = -50011966
$ws.Range("H126").Value = 637.8
$ws.Range("I126").Value = 655.625
$ws.Range("J126").Value = 566.5
$ws.Range("K126").Value = 1966.875
$ws.Range("L126").Value = 1699.5
$ws.Range("M126").Value = 503.125
$ws.Range("N126").Value = -6639.5
